$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "Darwin Nicolas,Mr. Howard Willms II"
    3  = "Osvaldo Boyle PhD"
    4  = "Carissa Adams"
    5  = "Carissa Adams"
    6  = "Prof. Eladio Franecki"
    7  = "Osvaldo Boyle PhD"
    8  = "Osvaldo Boyle PhD"
    9  = "Prof. Eladio Franecki"
    10 = "Darwin Nicolas"
    11 = "Dr. Yvonne Treutel PhD"
    12 = "Dr. Yvonne Treutel PhD"
    13 = "Prof. Eladio Franecki"
    14 = "Mr. Howard Willms II"
    15 = "Vernie Sporer"
    16 = "Dr. Yvonne Treutel PhD"
    17 = "Mr. Howard Willms II"
    18 = "Destinee Feest"
    19 = "Laisha Schultz"
    20 = "Carissa Adams"
    21 = "Laisha Schultz"
    22 = "Destinee Feest"
    23 = "Maymie Haley,Destinee Feest"
    24 = "Laisha Schultz"
    25 = "Maymie Haley"
    26 = "Maymie Haley"
    27 = "Mr. Howard Willms II"
    28 = "Destinee Feest,Vernie Sporer"
    29 = "Vernie Sporer"
    30 = "Darwin Nicolas"
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
}
